$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 221 (header + 220 data rows). Append
# three new data rows (222-224) with the same formatting as the rest of
# the table by copying the format of the last existing row down first.
$ws.Range("A221:E221").Copy()
$ws.Range("A222:E224").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 222
$ws.Cells.Item(222, 1).Value = 220
$ws.Cells.Item(222, 2).Value = 644
$ws.Cells.Item(222, 3).Value = 746
$ws.Cells.Item(222, 4).Value = 560
$ws.Cells.Item(222, 5).Value = "jc, marvel, lol, marbels"

# Row 223
$ws.Cells.Item(223, 1).Value = 221
$ws.Cells.Item(223, 2).Value = 622
$ws.Cells.Item(223, 3).Value = 925
$ws.Cells.Item(223, 4).Value = 270
$ws.Cells.Item(223, 5).Value = "jc, lol, cs"

# Row 224
$ws.Cells.Item(224, 1).Value = 222
$ws.Cells.Item(224, 2).Value = 916
$ws.Cells.Item(224, 3).Value = 1123
$ws.Cells.Item(224, 4).Value = 427
$ws.Cells.Item(224, 5).Value = "jc, cs, marbels"
